# "Natmi following Dr Hou advice"
#
# The natmi LR-pairs sheet (Tnfsf13 -> Tnfrsf11b) is recomputed: each
# sending cluster now yields a row for EVERY target cluster actually
# considered (ECs and FAPs), instead of a single collapsed row. This
# doubles the data rows (5 senders x 2 targets = 10) and refreshes all
# of the expression / specificity statistics accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf13"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.153770666666667
$ws.Range("H2").Value = 3.461312
$ws.Range("I2").Value = 0.1294297218267158
$ws.Range("J2").Value = 0.1294297218267158
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("N2").Value = 0.247237
$ws.Range("O2").Value = 0.04727005612861496
$ws.Range("P2").Value = 0.04727005612861496
$ws.Range("Q2").Value = 0.09508493277155555
$ws.Range("R2").Value = 0.8557643949439999
$ws.Range("S2").Value = 0.006118150215459878
$ws.Range("T2").Value = 0.006118150215459878

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf13"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.153770666666667
$ws.Range("H3").Value = 3.461312
$ws.Range("I3").Value = 0.1294297218267158
$ws.Range("J3").Value = 0.1294297218267158
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.661024
$ws.Range("N3").Value = 4.983072
$ws.Range("O3").Value = 0.9527299438713851
$ws.Range("P3").Value = 0.952729943871385
$ws.Range("Q3").Value = 1.916440767829333
$ws.Range("R3").Value = 17.247966910464
$ws.Range("S3").Value = 0.123311571611256
$ws.Range("T3").Value = 0.1233115716112559

# Row 4: FAPs -> ECs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf13"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.560947
$ws.Range("H4").Value = 1.682841
$ws.Range("I4").Value = 0.06292690243138796
$ws.Range("J4").Value = 0.06292690243138795
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("N4").Value = 0.247237
$ws.Range("O4").Value = 0.04727005612861496
$ws.Range("P4").Value = 0.04727005612861496
$ws.Range("Q4").Value = 0.04622895114633333
$ws.Range("R4").Value = 0.416060560317
$ws.Range("S4").Value = 0.002974558209931586
$ws.Range("T4").Value = 0.002974558209931586

# Row 5: FAPs -> FAPs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf13"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.560947
$ws.Range("H5").Value = 1.682841
$ws.Range("I5").Value = 0.06292690243138796
$ws.Range("J5").Value = 0.06292690243138795
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.661024
$ws.Range("N5").Value = 4.983072
$ws.Range("O5").Value = 0.9527299438713851
$ws.Range("P5").Value = 0.952729943871385
$ws.Range("Q5").Value = 0.931746429728
$ws.Range("R5").Value = 8.385717867552
$ws.Range("S5").Value = 0.05995234422145638
$ws.Range("T5").Value = 0.05995234422145636

# Row 6: M1 -> ECs
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Tnfsf13"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.837574333333334
$ws.Range("H6").Value = 11.512723
$ws.Range("I6").Value = 0.4304981854736105
$ws.Range("J6").Value = 0.4304981854736104
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("N6").Value = 0.247237
$ws.Range("O6").Value = 0.04727005612861496
$ws.Range("P6").Value = 0.04727005612861496
$ws.Range("Q6").Value = 0.3162634551501112
$ws.Range("R6").Value = 2.846371096351
$ws.Range("S6").Value = 0.02034967339060446
$ws.Range("T6").Value = 0.02034967339060446

# Row 7: M1 -> FAPs
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Tnfsf13"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.837574333333334
$ws.Range("H7").Value = 11.512723
$ws.Range("I7").Value = 0.4304981854736105
$ws.Range("J7").Value = 0.4304981854736104
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.661024
$ws.Range("N7").Value = 4.983072
$ws.Range("O7").Value = 0.9527299438713851
$ws.Range("P7").Value = 0.952729943871385
$ws.Range("Q7").Value = 6.374303069450668
$ws.Range("R7").Value = 57.368727625056
$ws.Range("S7").Value = 0.4101485120830061
$ws.Range("T7").Value = 0.4101485120830059

# Row 8: M2 -> ECs
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Tnfsf13"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.179188666666667
$ws.Range("H8").Value = 9.537566
$ws.Range("I8").Value = 0.3566406363494371
$ws.Range("J8").Value = 0.356640636349437
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("N8").Value = 0.247237
$ws.Range("O8").Value = 0.04727005612861496
$ws.Range("P8").Value = 0.04727005612861496
$ws.Range("Q8").Value = 0.2620043561268889
$ws.Range("R8").Value = 2.358039205142
$ws.Range("S8").Value = 0.01685842289798285
$ws.Range("T8").Value = 0.01685842289798285

# Row 9: M2 -> FAPs
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Tnfsf13"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.179188666666667
$ws.Range("H9").Value = 9.537566
$ws.Range("I9").Value = 0.3566406363494371
$ws.Range("J9").Value = 0.356640636349437
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.661024
$ws.Range("N9").Value = 4.983072
$ws.Range("O9").Value = 0.9527299438713851
$ws.Range("P9").Value = 0.952729943871385
$ws.Range("Q9").Value = 5.280708675861334
$ws.Range("R9").Value = 47.526378082752
$ws.Range("S9").Value = 0.3397822134514543
$ws.Range("T9").Value = 0.3397822134514542

# Row 10: sCs -> ECs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tnfsf13"
$ws.Range("C10").Value = "Tnfrsf11b"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.182783
$ws.Range("H10").Value = 0.548349
$ws.Range("I10").Value = 0.02050455391884863
$ws.Range("J10").Value = 0.02050455391884863
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("N10").Value = 0.247237
$ws.Range("O10").Value = 0.04727005612861496
$ws.Range("P10").Value = 0.04727005612861496
$ws.Range("Q10").Value = 0.01506357352366667
$ws.Range("R10").Value = 0.135572161713
$ws.Range("S10").Value = 0.0009692514146361869
$ws.Range("T10").Value = 0.0009692514146361867

# Row 11: sCs -> FAPs
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Tnfsf13"
$ws.Range("C11").Value = "Tnfrsf11b"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.182783
$ws.Range("H11").Value = 0.548349
$ws.Range("I11").Value = 0.02050455391884863
$ws.Range("J11").Value = 0.02050455391884863
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.661024
$ws.Range("N11").Value = 4.983072
$ws.Range("O11").Value = 0.9527299438713851
$ws.Range("P11").Value = 0.952729943871385
$ws.Range("Q11").Value = 0.303606949792
$ws.Range("R11").Value = 2.732462548128
$ws.Range("S11").Value = 0.01953530250421245
$ws.Range("T11").Value = 0.01953530250421244
